# kepsus.xlsx edit — IntensitasFilter: rename the "Rawan Bencana Banjir Tingkat
# Tinggi" label used by rows 188-322 (column B) to the new
# "Rawan Bencana Cuaca Ekstrem Tingkat Tinggi" label, and leave the sheet
# positioned/selected the way the author's session ended up (topLeftCell
# A181, active cell C189).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldLabel = "Rawan Bencana Banjir Tingkat Tinggi"
$newLabel = "Rawan Bencana Cuaca Ekstrem Tingkat Tinggi"

$firstRow = 188
$lastRow  = 322

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)   # column B
    if ($cell.Value2 -eq $oldLabel) {
        $cell.Value = $newLabel
    }
}

# Restore the view state captured in the saved workbook: scrolled so row 181
# is at the top, with C189 as the active/selected cell.
$ws.Activate()
$win = $excel.ActiveWindow
try { $win.ScrollRow = 181 } catch {}
try { $win.ScrollColumn = 1 } catch {}
$ws.Range("C189").Select()
